$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1541.6111
$ws.Range("I100").Value = 1174.9
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1174.9
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -633.9000000000001
$ws.Range("N100").Value = -3082

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1960.1
$ws.Range("I113").Value = 2099
$ws.Range("J113").Value = 1821.2
$ws.Range("K113").Value = 2099
$ws.Range("L113").Value = 1821.2
$ws.Range("M113").Value = 1155
$ws.Range("N113").Value = -8329.200000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 52926.668
$ws.Range("J133").Value = 52926.668
$ws.Range("L133").Value = 52926.668
$ws.Range("N133").Value = -63046.668

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1650
$ws.Range("I141").Value = 1650
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4950
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 230
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7356686
$ws.Range("I2").Value = 4993.3335
$ws.Range("K2").Value = 4993.3335
$ws.Range("M2").Value = -4880.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3264.4443
$ws.Range("I61").Value = 2803.4
$ws.Range("J61").Value = 3840.75
$ws.Range("K61").Value = 2803.4
$ws.Range("L61").Value = 3840.75
$ws.Range("M61").Value = -2591.4
$ws.Range("N61").Value = -4264.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1333.9231
$ws.Range("I110").Value = 1212.3334
$ws.Range("J110").Value = 1607.5
$ws.Range("K110").Value = 1212.3334
$ws.Range("L110").Value = 1607.5
$ws.Range("M110").Value = 832.6666
$ws.Range("N110").Value = -5697.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 7356686
$ws.Range("I116").Value = 4993.3335
$ws.Range("K116").Value = 4993.3335
$ws.Range("M116").Value = -2699.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2154.5
$ws.Range("I122").Value = 1962.2858
$ws.Range("K122").Value = 5886.857400000001
$ws.Range("M122").Value = -3436.857400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3264.4443
$ws.Range("I136").Value = 2803.4
$ws.Range("J136").Value = 3840.75
$ws.Range("K136").Value = 8410.200000000001
$ws.Range("L136").Value = 11522.25
$ws.Range("M136").Value = -5860.200000000001
$ws.Range("N136").Value = -16622.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7356686
$ws.Range("I3").Value = 4993.3335
$ws.Range("K3").Value = 4993.3335
$ws.Range("M3").Value = -4879.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 20666.334
$ws.Range("I26").Value = 20666.334
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 20666.334
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -20374.334
$ws.Range("N26").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1887.25
$ws.Range("I99").Value = 1966.3334
$ws.Range("J99").Value = 1650
$ws.Range("K99").Value = 1966.3334
$ws.Range("L99").Value = 1650
$ws.Range("M99").Value = -468.3334
$ws.Range("N99").Value = -4646

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2616.25
$ws.Range("I105").Value = 3320
$ws.Range("J105").Value = 1443.3334
$ws.Range("K105").Value = 3320
$ws.Range("L105").Value = 1443.3334
$ws.Range("M105").Value = -1573
$ws.Range("N105").Value = -4937.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1449.1724
$ws.Range("I107").Value = 1367.8096
$ws.Range("K107").Value = 1367.8096
$ws.Range("M107").Value = 552.1904

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 70540.44
$ws.Range("I134").Value = 70540.44
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 211621.32
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -209086.32
$ws.Range("N134").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8000
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 15000
$ws.Range("J28").Value = 15000
$ws.Range("L28").Value = 15000
$ws.Range("N28").Value = -15490

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3725.4119
$ws.Range("I132").Value = 3736
$ws.Range("K132").Value = 11208
$ws.Range("M132").Value = -8678

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3954.2354
$ws.Range("I134").Value = 4148.1333
$ws.Range("K134").Value = 12444.3999
$ws.Range("M134").Value = -9909.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 14286801
$ws.Range("I4").Value = 33334000
$ws.Range("J4").Value = 1402.5
$ws.Range("K4").Value = 100002000
$ws.Range("L4").Value = 4207.5
$ws.Range("M4").Value = -100001888
$ws.Range("N4").Value = -4431.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 166.83333
$ws.Range("I38").Value = 37
$ws.Range("J38").Value = 296.66666
$ws.Range("K38").Value = 111
$ws.Range("L38").Value = 889.9999799999999
$ws.Range("M38").Value = 236
$ws.Range("N38").Value = -1583.99998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 66668468
$ws.Range("I60").Value = 166666670
$ws.Range("J60").Value = 2993.3333
$ws.Range("K60").Value = 500000010
$ws.Range("L60").Value = 8979.999899999999
$ws.Range("M60").Value = -499999759
$ws.Range("N60").Value = -9481.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 588.0333000000001
$ws.Range("I113").Value = 621.2857
$ws.Range("K113").Value = 1863.8571
$ws.Range("M113").Value = 306.1428999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1726651.5
$ws.Range("I131").Value = 4901.8184
$ws.Range("K131").Value = 14705.4552
$ws.Range("M131").Value = -9665.4552

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 41.411766
$ws.Range("I2").Value = 33.545456
$ws.Range("J2").Value = 55.833332
$ws.Range("K2").Value = 33.545456
$ws.Range("L2").Value = 55.833332
$ws.Range("M2").Value = 79.454544
$ws.Range("N2").Value = -281.833332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2192.397
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 2192.397
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 2192.397
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = -2416.397

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3476.1667
$ws.Range("I122").Value = 3521.4
$ws.Range("K122").Value = 10564.2
$ws.Range("M122").Value = -8114.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2387562.8
$ws.Range("I2").Value = 250030
$ws.Range("J2").Value = 3337577.5
$ws.Range("K2").Value = 250030
$ws.Range("L2").Value = 3337577.5
$ws.Range("M2").Value = -249918
$ws.Range("N2").Value = -3337801.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1971.1818
$ws.Range("I7").Value = 1867.5555
$ws.Range("J7").Value = 2437.5
$ws.Range("K7").Value = 1867.5555
$ws.Range("L7").Value = 2437.5
$ws.Range("M7").Value = -1755.5555
$ws.Range("N7").Value = -2661.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1971.1818
$ws.Range("I126").Value = 1867.5555
$ws.Range("J126").Value = 2437.5
$ws.Range("K126").Value = 5602.666499999999
$ws.Range("L126").Value = 7312.5
$ws.Range("M126").Value = -3132.666499999999
$ws.Range("N126").Value = -12252.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8004.7144
$ws.Range("I132").Value = 10677.385
$ws.Range("J132").Value = 3661.625
$ws.Range("K132").Value = 32032.155
$ws.Range("L132").Value = 10984.875
$ws.Range("M132").Value = -29502.155
$ws.Range("N132").Value = -16044.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 15026500
$ws.Range("I2").Value = 2038400
$ws.Range("J2").Value = 36673336
$ws.Range("K2").Value = 2038400
$ws.Range("L2").Value = 36673336
$ws.Range("M2").Value = -2038288
$ws.Range("N2").Value = -36673560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1963.2858
$ws.Range("I132").Value = 842.55554
$ws.Range("J132").Value = 3980.6
$ws.Range("K132").Value = 2527.66662
$ws.Range("L132").Value = 11941.8
$ws.Range("M132").Value = 2.333380000000034
$ws.Range("N132").Value = -17001.8
